$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.232.63"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "1.851.61"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4646"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3712"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07289"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8870"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07850"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").Value = "1.839.82"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.399"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.509"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008921"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "27.256.19"
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("D24").Value = "2.097.69"
$ws.Range("E24").Value = "  +3.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.955"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.040"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.99"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.032"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08835"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.143"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7659"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.19%  "

$ws.Range("E34").Value = "  +2.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.521"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.60%  "

$ws.Range("E36").Value = "  +10.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01940"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05218"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.946"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("E41").Value = "  -1.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1625"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.468"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4794"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.642"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06205"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
